$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B25 becomes a numeric value instead of inline string "2"
$ws.Range("B25").Value = 2

# Add new row 26 with data
$ws.Range("A26").Value = "Ying Tang"
$ws.Range("B26").Value = "4"
$ws.Range("C26").Value = "It took us as a few weeks to reply,we took the time to implement ,as much as possible"
$ws.Range("D26").Value = "ACK"
$ws.Range("E26").Value = "OTH"
$ws.Range("F26").Value = "91b1b71f-4957-400a-bdb5-bced2ed448de"
$ws.Range("G26").Value = "S1CChZ-CZ_annotated.xlsx"
$ws.Range("H26").Value = "It took us as a few weeks to reply because we took the time to implement as much as possible of the feedback."
